$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Grupo 5 (row 8): "William Omar Cueto" -> "Matias Villalba" (name correction, no accent)
$ws.Range("C8").Value = "Matias Villalba"

# 2. Delete the whole row for old Grupo 14 (Sindy Gonzalez / Matias Villalba) - row 17.
#    Everything below shifts up by one row.
$ws.Rows("17:17").Delete()

# 3. Add a note in column E next to Grupo 13 (now row 16, Agustin Giuli / Lorenzo Vera)
$ws.Range("E16").Value = "   Agustin Giuli participo ?"

# 4. The two members who lost their partner (William Omar Cueto and Sindy Gonzalez)
#    are re-listed standalone as new Grupo 18 and Grupo 19 (rows 20 and 21),
#    re-using the formatting of the other "name" cells (like C6).
$ws.Range("C6").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = "William Omar Cueto"

$ws.Range("C6").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C21").Value = "Sindy González"

$ws.Range("A1").Select()
